$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.172316384180791
$ws.Range("C2").Value = 0.576271186440678
$ws.Range("J2").Value = 0.01694915254237288
$ws.Range("P2").Value = 0.1299435028248588
$ws.Range("S2").Value = 0.1045197740112994
$ws.Range("B3").Value = 0.009433962264150943
$ws.Range("C3").Value = 0.01886792452830189
$ws.Range("J3").Value = 0.04716981132075472
$ws.Range("P3").Value = 0.7216981132075472
$ws.Range("S3").Value = 0.2028301886792453
$ws.Range("J4").Value = 0.02083333333333333
$ws.Range("P4").Value = 0.7083333333333334
$ws.Range("S4").Value = 0.2708333333333333
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.05527638190954774
$ws.Range("D6").Value = 0.01507537688442211
$ws.Range("E6").Value = 0.005025125628140704
$ws.Range("F6").Value = 0.03015075376884422
$ws.Range("J6").Value = 0.3115577889447236
$ws.Range("O6").Value = 0.01507537688442211
$ws.Range("Q6").Value = 0.1959798994974874
$ws.Range("R6").Value = 0.03517587939698492
$ws.Range("S6").Value = 0.3366834170854272
$ws.Range("B7").Value = 0.1151515151515152
$ws.Range("D7").Value = 0.01818181818181818
$ws.Range("F7").Value = 0.0303030303030303
$ws.Range("J7").Value = 0.1878787878787879
$ws.Range("O7").Value = 0.01818181818181818
$ws.Range("Q7").Value = 0.2181818181818182
$ws.Range("R7").Value = 0.03636363636363636
$ws.Range("S7").Value = 0.3757575757575757
$ws.Range("B8").Value = 0.131578947368421
$ws.Range("D8").Value = 0.02105263157894737
$ws.Range("E8").Value = 0.002631578947368421
$ws.Range("F8").Value = 0.05789473684210526
$ws.Range("J8").Value = 0.1026315789473684
$ws.Range("O8").Value = 0.01842105263157895
$ws.Range("Q8").Value = 0.1868421052631579
$ws.Range("R8").Value = 0.09210526315789473
$ws.Range("S8").Value = 0.3868421052631579
$ws.Range("B9").Value = 0.1372549019607843
$ws.Range("D9").Value = 0.006535947712418301
$ws.Range("E9").Value = 0.006535947712418301
$ws.Range("F9").Value = 0.0392156862745098
$ws.Range("J9").Value = 0.1372549019607843
$ws.Range("O9").Value = 0.0196078431372549
$ws.Range("Q9").Value = 0.2026143790849673
$ws.Range("R9").Value = 0.0915032679738562
$ws.Range("S9").Value = 0.3594771241830065
$ws.Range("B10").Value = 0.1348396501457726
$ws.Range("D10").Value = 0.02551020408163265
$ws.Range("F10").Value = 0.0641399416909621
$ws.Range("J10").Value = 0.130466472303207
$ws.Range("O10").Value = 0.02113702623906705
$ws.Range("Q10").Value = 0.2084548104956268
$ws.Range("R10").Value = 0.06851311953352769
$ws.Range("S10").Value = 0.3469387755102041
$ws.Range("G11").Value = 0.1283018867924528
$ws.Range("J11").Value = 0.1169811320754717
$ws.Range("K11").Value = 0.1773584905660377
$ws.Range("L11").Value = 0.5735849056603773
$ws.Range("S11").Value = 0.003773584905660377
$ws.Range("G12").Value = 0.7161290322580646
$ws.Range("J12").Value = 0.2387096774193548
$ws.Range("K12").Value = 0.006451612903225806
$ws.Range("L12").Value = 0.01290322580645161
$ws.Range("S12").Value = 0.02580645161290323
$ws.Range("G13").Value = 0.575
$ws.Range("J13").Value = 0.4
$ws.Range("S13").Value = 0.025
$ws.Range("F15").Value = 0.01953125
$ws.Range("H15").Value = 0.1328125
$ws.Range("I15").Value = 0.078125
$ws.Range("J15").Value = 0.40234375
$ws.Range("K15").Value = 0.0859375
$ws.Range("M15").Value = 0.00390625
$ws.Range("O15").Value = 0.0625
$ws.Range("S15").Value = 0.21484375
$ws.Range("F16").Value = 0.01339285714285714
$ws.Range("H16").Value = 0.1517857142857143
$ws.Range("I16").Value = 0.09821428571428571
$ws.Range("J16").Value = 0.3973214285714285
$ws.Range("K16").Value = 0.1071428571428571
$ws.Range("M16").Value = 0.03125
$ws.Range("O16").Value = 0.07589285714285714
$ws.Range("S16").Value = 0.125
$ws.Range("F17").Value = 0.02608695652173913
$ws.Range("H17").Value = 0.1847826086956522
$ws.Range("I17").Value = 0.05434782608695652
$ws.Range("J17").Value = 0.4565217391304348
$ws.Range("K17").Value = 0.08043478260869565
$ws.Range("M17").Value = 0.02391304347826087
$ws.Range("O17").Value = 0.06521739130434782
$ws.Range("S17").Value = 0.108695652173913
$ws.Range("F18").Value = 0.01298701298701299
$ws.Range("H18").Value = 0.1558441558441558
$ws.Range("I18").Value = 0.09090909090909091
$ws.Range("J18").Value = 0.474025974025974
$ws.Range("K18").Value = 0.08441558441558442
$ws.Range("M18").Value = 0.006493506493506494
$ws.Range("O18").Value = 0.07792207792207792
$ws.Range("S18").Value = 0.09740259740259741
$ws.Range("F19").Value = 0.02480752780153978
$ws.Range("H19").Value = 0.1753635585970915
$ws.Range("I19").Value = 0.06330196749358426
$ws.Range("J19").Value = 0.4080410607356715
$ws.Range("K19").Value = 0.1043627031650984
$ws.Range("M19").Value = 0.01967493584260051
$ws.Range("O19").Value = 0.08982035928143713
$ws.Range("S19").Value = 0.1146278870829769
